$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update loading_percent values for rows 2-25 (indices 0-23), columns B:N
# Columns F, H, J remain 0 (unchanged); only B,C,D,E,G,I,K,L,M,N change per row.
$data = @{
    2 = @(18.98675272900141, 8.189968631619857, 6.077014308638037, 8.514233786460656, 0, 3.775282611885188, 0, 38.26958236946545, 0, 17.47447582109373, 10.61722788461962, 17.50285539737063, 26.07614836210371)
    3 = @(18.85368605590368, 7.976876354802147, 5.973751394429807, 8.523139282333672, 0, 3.778799528406593, 0, 38.21137738503783, 0, 17.36959470919817, 10.62794495146161, 17.4984062799528, 26.0836894591746)
    4 = @(18.77680316077571, 7.845738182111106, 5.911349246916056, 8.52902909165541, 0, 3.781071197545371, 0, 38.18105643295558, 0, 17.30949878069627, 10.63603937015721, 17.49902539142119, 26.09002748644447)
    5 = @(18.74671383026947, 7.792311286409932, 5.886207405695036, 8.531535527392304, 0, 3.782025256173283, 0, 38.17006675847568, 0, 17.28611119229516, 10.63971876546717, 17.50012144023679, 26.09303961675573)
    6 = @(18.74179327707895, 7.783443210067266, 5.882051023602099, 8.531958145514958, 0, 3.782185391310113, 0, 38.16832456175819, 0, 17.28229481666758, 10.64035273121938, 17.50035442460914, 26.09356570567191)
    7 = @(18.77639230405446, 7.845017477058835, 5.911008964141931, 8.52906246364204, 0, 3.781083949449534, 0, 38.18090268526949, 0, 17.3091788801248, 10.63608744959586, 17.49903675560054, 26.09006637097569)
    8 = @(18.9398888542703, 8.116608198705846, 6.041222829257112, 8.517216992210392, 0, 3.776472006323913, 0, 38.24839169722385, 0, 17.43743279728672, 10.62060896850443, 17.5006266693964, 26.07839403345078)
    9 = @(19.29737905323805, 8.643142283625629, 6.302945979490208, 8.497324920968975, 0, 3.768314041759806, 0, 38.42350816709328, 0, 17.72209032387073, 10.60226264471683, 17.53026685652657, 26.06906032913123)
    10 = @(19.58052855529076, 9.021669642307973, 6.497028536767462, 8.484731121269363, 0, 3.762853884900178, 0, 38.57791522223158, 0, 17.95001574683416, 10.59609218316604, 17.56810242645104, 26.07047558944416)
    11 = @(19.71333681280348, 9.191185665739763, 6.58528100819562, 8.479437925499345, 0, 3.760484342097333, 0, 38.6536640117954, 0, 18.05745134289961, 10.59486836085039, 17.58876976893945, 26.07291705348422)
    12 = @(19.76416317835761, 9.254925806290993, 6.618660193892153, 8.477495981524905, 0, 3.759603386643789, 0, 38.68313130532199, 0, 18.09864398849041, 10.59463211986537, 17.59708938940537, 26.07410003338939)
    13 = @(19.75319368591305, 9.241219330568148, 6.611473712958668, 8.477911438355182, 0, 3.759792391047524, 0, 38.67675033334015, 0, 18.08975026463068, 10.59467290224548, 17.59527573029951, 26.07383376411893)
    14 = @(19.7175078477249, 9.196439025437387, 6.58802813241487, 8.479276909602964, 0, 3.760411538467749, 0, 38.65607266953776, 0, 18.06083026012331, 10.59484437477459, 17.58944436381755, 26.07300919867781)
    15 = @(19.69571765797806, 9.168949055474085, 6.573660803246642, 8.480121430561253, 0, 3.760792909037732, 0, 38.64350866896968, 0, 18.04318130871219, 10.59497897899855, 17.58593661785208, 26.07253778412196)
    16 = @(19.57192629599249, 9.010531648068794, 6.491257283829325, 8.485085796786546, 0, 3.76301103185972, 0, 38.57307493237905, 0, 17.94306749395939, 10.59620398222193, 17.56682100896429, 26.07035220421077)
    17 = @(19.49698067808, 8.912612538114187, 6.440670798167646, 8.488242753041675, 0, 3.764400985535883, 0, 38.53127062136356, 0, 17.88259044285002, 10.59736066139358, 17.5559768900369, 26.06947182326418)
    18 = @(19.45425314937079, 8.856043641226, 6.411573344993393, 8.490099580340365, 0, 3.765211215390655, 0, 38.50774495375161, 0, 17.84816060390043, 10.59817499889449, 17.550065153957, 26.06913466122708)
    19 = @(19.43985269681709, 8.836849844485741, 6.401722282541712, 8.490735322907396, 0, 3.765487397254501, 0, 38.49986898401369, 0, 17.8365651236401, 10.59847633119096, 17.54811955117821, 26.06904956662546)
    20 = @(19.50491981257701, 8.923062405798188, 6.446056222052627, 8.487902444635965, 0, 3.764251909150173, 0, 38.5356671045685, 0, 17.88899180747297, 10.5972221083211, 17.55709760018484, 26.06954802910141)
    21 = @(19.72797545944124, 9.209604856291975, 6.594916027151452, 8.47887414356955, 0, 3.76022923721056, 0, 38.66212503118602, 0, 18.06931119538671, 10.59478784721741, 17.59114381896818, 26.07324437988992)
    22 = @(19.8768527208113, 9.394215758385643, 6.691955569904925, 8.473337689044671, 0, 3.757695379194393, 0, 38.74933092131698, 0, 18.1901131739701, 10.59452090836009, 17.61626862176022, 26.07716654975896)
    23 = @(19.7971243881466, 9.29594960332231, 6.640197664551637, 8.476259349640916, 0, 3.759039069424234, 0, 38.70237366244106, 0, 18.12537894741627, 10.59454240576121, 17.60259740128515, 26.07493540755424)
    24 = @(19.50132940633901, 8.918338868630796, 6.443621514779798, 8.488056167703405, 0, 3.764319271917031, 0, 38.53367787080305, 0, 17.88609669159167, 10.5972842829372, 17.5565899217079, 26.06951305002993)
    25 = @(19.19691808343973, 8.501845337279162, 6.231678539555376, 8.502350409505748, 0, 3.770426820234693, 0, 38.37158580532039, 0, 17.64167462413421, 10.60594125565759, 17.5194186204304, 26.07013302974105)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $arr = New-Object "object[,]" 1,13
    for ($i = 0; $i -lt 13; $i++) {
        $arr[0, $i] = $vals[$i]
    }
    $ws.Range("B$row`:N$row").Value = $arr
}
